# Add "minimum_op_point_Output1" / "minimum_op_point_Output2" to the Units
# table (Table1) and populate the Electrolyzer row's ramp/shut-down/minimum
# operating point figures, matching the "Added minimum operating point to
# spine and jupyter" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Units")
$ws.Activate()

$lo = $ws.ListObjects.Item("Table1")

# --- Extend the table with the two new trailing columns -------------------
$col1 = $lo.ListColumns.Add()
$ws.Cells.Item(1, $col1.Index).Value = "minimum_op_point_Output1"

$col2 = $lo.ListColumns.Add()
$ws.Cells.Item(1, $col2.Index).Value = "minimum_op_point_Output2"

# --- Update the Solar_Plant_Kasso (row 2) figures --------------------------
$ws.Range("O2").Value = 0.3    # ramp_up_Output1
$ws.Range("Q2").Value = 0.1    # ramp_down_Output1
$ws.Range("U2").Value = 0.2    # shut_down_Output1 (was 0.8)
$ws.Range("AH2").Value = 0.2   # minimum_op_point_Output1

# --- Match the author's last on-screen selection ---------------------------
$ws.Range("AH1").Select()
